# Actualización automática 2025-11-10 16:30:09
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Adjust column widths for E (5) and F (6)
# (ColumnWidth has a constant +5/6 offset baked into the resulting OOXML
#  "width" attribute, so we back it out here to land on the exact targets)
$ws.Columns.Item(5).ColumnWidth = 21.166666666666668
$ws.Columns.Item(6).ColumnWidth = 25.166666666666668

# Row 2 - 240X120 PORCELANATO
$ws.Cells.Item(2, 3).Value = 129.6
$ws.Cells.Item(2, 5).Value = 129.6

# Row 3 - 240X80 PORCELANATO
$ws.Cells.Item(3, 3).Value = 2564
$ws.Cells.Item(3, 4).Value = -152.64
$ws.Cells.Item(3, 5).Value = 2716.64
$ws.Cells.Item(3, 6).Value = -0.05953198127925116

# Row 4 - FREGADEROS DE COCINA
$ws.Cells.Item(4, 3).Value = 207.39
$ws.Cells.Item(4, 5).Value = 207.39

# Row 5 - GRIFERIAS
$ws.Cells.Item(5, 3).Value = 86.41
$ws.Cells.Item(5, 5).Value = 86.41

# Row 6 - INODOROS
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 2907.58368146026
$ws.Cells.Item(6, 6).Value = 0

# Row 7 - LAVABOS
$ws.Cells.Item(7, 3).Value = 383.4
$ws.Cells.Item(7, 5).Value = 383.4

# Row 8 - NO RESURTIBLES
$ws.Cells.Item(8, 3).Value = 415
$ws.Cells.Item(8, 5).Value = 415

# Row 10 - PANELES DECORATIVOS
$ws.Cells.Item(10, 3).Value = 1388
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 1388
$ws.Cells.Item(10, 6).Value = 0

# Row 11 - PIEDRA SINTERIZADA
$ws.Cells.Item(11, 3).Value = 2678
$ws.Cells.Item(11, 4).Value = -1151.4
$ws.Cells.Item(11, 5).Value = 3829.4
$ws.Cells.Item(11, 6).Value = -0.4299477221807319

# Row 12 - PORCELANATO
$ws.Cells.Item(12, 3).Value = 44418
$ws.Cells.Item(12, 4).Value = 3045.53
$ws.Cells.Item(12, 5).Value = 41372.47
$ws.Cells.Item(12, 6).Value = 0.06856522130667747

# Row 14 - TOTAL
$ws.Cells.Item(14, 3).Value = 55399.47101170094
$ws.Cells.Item(14, 4).Value = 1741.49
$ws.Cells.Item(14, 5).Value = 53657.98101170095
$ws.Cells.Item(14, 6).Value = 0.03143513770433259
